$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 5934
$ws.Range("K3").Value = 6104
$ws.Range("K4").Value = 1272
$ws.Range("I5").Value = 576
$ws.Range("K5").Value = 433
$ws.Range("K6").Value = 6712
$ws.Range("I7").Value = 20780
$ws.Range("K7").Value = 20455

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("K6").Value = 11
$ws.Range("K7").Value = 46

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 375
$ws.Range("K3").Value = 415
$ws.Range("K6").Value = 450
$ws.Range("K7").Value = 1353

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 234
$ws.Range("K3").Value = 325
$ws.Range("K7").Value = 884

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 118
$ws.Range("K7").Value = 345

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K3").Value = 229
$ws.Range("K7").Value = 692

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 158
$ws.Range("K7").Value = 478

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K5").Value = 10
$ws.Range("K7").Value = 336

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K5").Value = 51
$ws.Range("K7").Value = 595
$ws.Range("K8").Value = 1353
$ws.Range("K10").Value = 114
$ws.Range("K11").Value = 386
$ws.Range("K13").Value = 28
$ws.Range("K14").Value = 105
$ws.Range("K19").Value = 587
$ws.Range("K20").Value = 483
$ws.Range("K23").Value = 210
$ws.Range("K24").Value = 60
$ws.Range("K26").Value = 26
$ws.Range("K28").Value = 8
$ws.Range("K29").Value = 1119
$ws.Range("K33").Value = 884
$ws.Range("K34").Value = 116
$ws.Range("K36").Value = 261
$ws.Range("K37").Value = 692
$ws.Range("K42").Value = 762
$ws.Range("K45").Value = 27
$ws.Range("K48").Value = 256
$ws.Range("K50").Value = 98
$ws.Range("K51").Value = 263
$ws.Range("K52").Value = 541
$ws.Range("K54").Value = 398
$ws.Range("K55").Value = 228
$ws.Range("K59").Value = 34
$ws.Range("I63").Value = 185
$ws.Range("K63").Value = 58
$ws.Range("K65").Value = 478
$ws.Range("K67").Value = 801
$ws.Range("K69").Value = 46
$ws.Range("K78").Value = 232
$ws.Range("K79").Value = 506
$ws.Range("K80").Value = 72
$ws.Range("K86").Value = 129
$ws.Range("K89").Value = 300
$ws.Range("K90").Value = 189
$ws.Range("K91").Value = 231
$ws.Range("K93").Value = 75
$ws.Range("K95").Value = 345
$ws.Range("K99").Value = 336
$ws.Range("K100").Value = 39
$ws.Range("I101").Value = 20780
$ws.Range("K101").Value = 20455

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K3").Value = 284
$ws.Range("K7").Value = 801

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K2").Value = 65
$ws.Range("K6").Value = 215
$ws.Range("K7").Value = 398

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 318
$ws.Range("K3").Value = 402
$ws.Range("K4").Value = 54
$ws.Range("K6").Value = 317
$ws.Range("K7").Value = 1119

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K6").Value = 123
$ws.Range("K7").Value = 256

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 180
$ws.Range("K7").Value = 587

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 105

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 207
$ws.Range("K4").Value = 32
$ws.Range("K6").Value = 284
$ws.Range("K7").Value = 762

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("K5").Value = 10
$ws.Range("K6").Value = 28

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K2").Value = 34
$ws.Range("K7").Value = 114

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K2").Value = 69
$ws.Range("K3").Value = 56
$ws.Range("K7").Value = 232

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K2").Value = 71
$ws.Range("K3").Value = 63
$ws.Range("K7").Value = 228

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("K3").Value = 13
$ws.Range("K7").Value = 60

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K2").Value = 60
$ws.Range("K7").Value = 210

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K4").Value = 9
$ws.Range("K7").Value = 231

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 171
$ws.Range("K7").Value = 506

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 161
$ws.Range("K3").Value = 157
$ws.Range("K4").Value = 21
$ws.Range("K7").Value = 483

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 103
$ws.Range("K3").Value = 77
$ws.Range("K7").Value = 261

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("K6").Value = 28
$ws.Range("K7").Value = 75

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("K6").Value = 22
$ws.Range("K7").Value = 39

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 198
$ws.Range("K3").Value = 194
$ws.Range("K7").Value = 595

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K2").Value = 45
$ws.Range("K7").Value = 116

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K2").Value = 26
$ws.Range("K7").Value = 98

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("K6").Value = 19
$ws.Range("K7").Value = 26

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K6").Value = 124
$ws.Range("K7").Value = 386

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K3").Value = 46
$ws.Range("K4").Value = 13

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("K2").Value = 9
$ws.Range("K7").Value = 34

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K2").Value = 84
$ws.Range("K6").Value = 91
$ws.Range("K7").Value = 300

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("K6").Value = 24
$ws.Range("K7").Value = 51

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 129

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K6").Value = 45
$ws.Range("K7").Value = 189

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K3").Value = 70
$ws.Range("K4").Value = 30
$ws.Range("K6").Value = 84
$ws.Range("K7").Value = 263

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("K2").Value = 6
$ws.Range("K7").Value = 27

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 72

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 146
$ws.Range("K3").Value = 157
$ws.Range("K6").Value = 190
$ws.Range("K7").Value = 541

$ws = $wb.Worksheets.Item("Edison Park")
$ws.Range("K6").Value = 1
$ws.Range("K7").Value = 8
